$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values (some values change)
$ws.Range("A2").Value = "13-04-2023"
$ws.Range("C2").Value = "10/04/2023  09:19:30"
$ws.Range("D2").Value = 12161100
$ws.Range("E2").Value = 8597675
$ws.Range("F2").Value = "CAIO HENRIQUE RODRIGUES FERNANDES"

# Add new row 3, copy of row 2 pattern but with its own values
$ws.Range("A3").Value = "13-04-2023"
$ws.Range("B3").Value = "417823 - PREMIUM SAÚDE S.A"
$ws.Range("C3").Value = "10/04/2023  10:47:50"
$ws.Range("D3").Value = 12161386
$ws.Range("E3").Value = 8598016
$ws.Range("F3").Value = "DIEGO SANTOS DE ALMEIDA"
$ws.Range("G3").Value = "8 dias úteis"
$ws.Range("H3").Value = "NO"
$ws.Range("I3").Value = "Assistencial"
$ws.Range("J3").Value = "Responder  Detalhes"
